$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Maximum Length" column (D) of data, added alongside the existing
# USB / SATA / Thunderbolt / WiFi throughput table.
$ws.Range("D1").Value = "Maximum Length"
$ws.Range("D3").Value = "4.5m"
$ws.Range("D6").Value = "1m"
$ws.Range("D10").Value = "3m"
$ws.Range("D14").Value = "5m"

# Footnote text at the bottom, now mentioning distances too.
$ws.Range("D23").Value = "*Theoretical throughputs and distances "

# WiFi cable-length column entries.
$ws.Range("D18").Value = "35m"
$ws.Range("D19").Value = "35m"
$ws.Range("D20").Value = "35m"
$ws.Range("D21").Value = "75m"
$ws.Range("D22").Value = "35m"

# Widen the new column to fit its content.
$ws.Range("D1").ColumnWidth = 20.25

# Fit the printed sheet to one page.
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1

# Leave the cursor where the author left it.
$ws.Range("B9").Select() | Out-Null
